# Update the "table_of_functions.csv" sheet to add the new ggflow function
# and mark pplot as deprecated (commit: "updated function list to include ggflow")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_of_functions.csv")

# Row 15 currently documents **pplot** -> "Draw a FlowJo-style (sort of) Plot"
# Mark it as deprecated.
$ws.Range("C15").Value = "Draw a FlowJo-style (sort of) Plot (DEPRECATED)"

# Insert a new row right after it (row 16) for the new, more powerful **ggflow** function,
# shifting everything below down by one row.
$ws.Rows.Item(16).Insert()
$ws.Range("B16").Value = "**ggflow**"
$ws.Range("C16").Value = "Draw a FlowJo-style (sort of) Plot (much more powerful version of pplot)"
